$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates in header (shared-string rich-text cells) ---
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# --- Crime statistics table updates (rows 14-27) ---
$ws.Range("N14").Value = -75
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 3
$ws.Range("C23").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 62.5
$ws.Range("N15").Value = -23.529411764705
$ws.Range("C16").Value = 7
$ws.Range("C38").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 3
$ws.Range("K38").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = 133.333333333333
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 106
$ws.Range("J16").Value = 81
$ws.Range("K16").Value = 30.864197530864
$ws.Range("L16").Value = 41.333333333333
$ws.Range("M16").Value = -15.873015873015
$ws.Range("N16").Value = -80.935251798561
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 154
$ws.Range("J17").Value = 123
$ws.Range("K17").Value = 25.203252032520
$ws.Range("L17").Value = 18.461538461538
$ws.Range("M17").Value = 87.804878048780
$ws.Range("N17").Value = -28.703703703703
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 13.333333333333
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 131
$ws.Range("K18").Value = 22.137404580152
$ws.Range("L18").Value = 64.948453608247
$ws.Range("M18").Value = -21.182266009852
$ws.Range("N18").Value = -85.388127853881
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -31.25
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 75
$ws.Range("H19").Value = -41.333333333333
$ws.Range("I19").Value = 472
$ws.Range("J19").Value = 531
$ws.Range("K19").Value = -11.111111111111
$ws.Range("L19").Value = 18.592964824120
$ws.Range("M19").Value = 60.544217687074
$ws.Range("N19").Value = -13.235294117647
$ws.Range("C20").Value = 3
$ws.Range("C38").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = 5
$ws.Range("K38").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = -40
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 38.461538461538
$ws.Range("I20").Value = 116
$ws.Range("J20").Value = 104
$ws.Range("K20").Value = 11.538461538461
$ws.Range("L20").Value = 84.126984126984
$ws.Range("M20").Value = -2.521008403361
$ws.Range("N20").Value = -90.951638065522
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -3.333333333333
$ws.Range("F21").Value = 114
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = -7.317073170731
$ws.Range("I21").Value = 1023
$ws.Range("J21").Value = 988
$ws.Range("K21").Value = 3.542510121457
$ws.Range("L21").Value = 31.660231660231
$ws.Range("M21").Value = 22.076372315035
$ws.Range("N21").Value = -72.485207100591
$ws.Range("L22").Value = 0
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = -49.019607843137
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 175
$ws.Range("H24").Value = -30.857142857142
$ws.Range("I24").Value = 1286
$ws.Range("J24").Value = 1316
$ws.Range("K24").Value = -2.279635258358
$ws.Range("L24").Value = 49.361207897793
$ws.Range("M24").Value = 82.930298719772
$ws.Range("C25").Value = 20
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 45.454545454545
$ws.Range("I25").Value = 423
$ws.Range("J25").Value = 320
$ws.Range("K25").Value = 32.1875
$ws.Range("L25").Value = 36.893203883495
$ws.Range("M25").Value = 33.438485804416
$ws.Range("C26").Value = 1
$ws.Range("C23").Copy() | Out-Null
$ws.Range("G26").PasteSpecial(-4122) | Out-Null
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "0"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("G26").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "***.*"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = -5.555555555555
$ws.Range("L26").Value = 13.333333333333
$ws.Range("C38").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 44
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = 18.918918918918
$ws.Range("L27").Value = 7.317073170731

$excel.CutCopyMode = 0
